$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1876
$ws1.Range("F7").Value = 3773
$ws1.Range("F8").Value = 183
$ws1.Range("F10").Value = 91
$ws1.Range("F12").Value = 82
$ws1.Range("F14").Value = 162
$ws1.Range("F15").Value = 900
$ws1.Range("F18").Value = 150
$ws1.Range("F20").Value = 93
$ws1.Range("F22").Value = 3264
$ws1.Range("F23").Value = 5607
$ws1.Range("F25").Value = 13
$ws1.Range("F30").Value = 342
$ws1.Range("F35").Value = 179
$ws1.Range("F37").Value = 341
$ws1.Range("F38").Value = 99
$ws1.Range("F39").Value = 494
$ws1.Range("F40").Value = 872
$ws1.Range("F41").Value = 43
$ws1.Range("F45").Value = 534

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 87

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1876
$ws4.Range("F7").Value = 3773
$ws4.Range("F8").Value = 183
$ws4.Range("F10").Value = 91
$ws4.Range("F11").Value = 87
$ws4.Range("F13").Value = 82
$ws4.Range("F15").Value = 162
$ws4.Range("F16").Value = 900
$ws4.Range("F19").Value = 150
$ws4.Range("F21").Value = 93
$ws4.Range("F23").Value = 3264
$ws4.Range("F24").Value = 5607
$ws4.Range("F26").Value = 13
$ws4.Range("F31").Value = 342
$ws4.Range("F36").Value = 179
$ws4.Range("F38").Value = 341
$ws4.Range("F39").Value = 99
$ws4.Range("F40").Value = 494
$ws4.Range("F41").Value = 872
$ws4.Range("F42").Value = 43
$ws4.Range("F46").Value = 534
